$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1 (05:22 -> 05:52)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 05:52"

# Row 19: India - updated case numbers
$ws.Cells.Item(19, 2).Value = 26496   # Casos totales
$ws.Cells.Item(19, 3).Value = 213     # Nuevos casos
$ws.Cells.Item(19, 5).Value = 19732   # Recuperados

# Row 125: El Salvador - updated case numbers
$ws.Cells.Item(125, 2).Value = 298    # Casos totales
$ws.Cells.Item(125, 3).Value = 24     # Nuevos casos
$ws.Cells.Item(125, 4).Value = 83     # Casos activos
$ws.Cells.Item(125, 5).Value = 207    # Recuperados

# Rows 173/174: Mongolia and San Martin (Parte Francesa) swap order
# and Mongolia's data is updated for the new day.
# Row 173 becomes Mongolia (updated numbers)
$ws.Cells.Item(173, 1).Value = "Mongolia"
$ws.Cells.Item(173, 2).Value = 38
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(173, 4).Value = 9
$ws.Cells.Item(173, 5).Value = 29
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

# Row 174 becomes San Martin (Parte Francesa) (unchanged numbers, just moved down)
$ws.Cells.Item(174, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(174, 2).Value = 38
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 24
$ws.Cells.Item(174, 5).Value = 11
$ws.Cells.Item(174, 6).Value = 3
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 3
